# Follow Map Time form: the "date" field (FMT_FOL_date) and the "time"
# field (FMT_time) are no longer backed by dedicated date/time widgets -
# both are now stored as plain "text" columns (per commit message:
# "Changed date and time fields for JGI app and verified database
# persistence").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Row 2 = FMT_FOL_date ("type" column was "date")
$ws.Range("C2").Value = "text"

# Row 4 = FMT_time ("type" column was "time")
$ws.Range("C4").Value = "text"

# Update the sheet's last-selected cell to match the saved view state.
$ws.Range("C5").Select()
